$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 662.75
$ws.Range("I2").Value = 454.7857
$ws.Range("J2").Value = 1148
$ws.Range("K2").Value = 454.7857
$ws.Range("L2").Value = 1148
$ws.Range("M2").Value = -341.7857
$ws.Range("N2").Value = -1374
$ws.Range("H19").Value = 1871.2
$ws.Range("I19").Value = 913.6
$ws.Range("K19").Value = 913.6
$ws.Range("M19").Value = -738.6
$ws.Range("H70").Value = 4916.6895
$ws.Range("I70").Value = 1217.6666
$ws.Range("J70").Value = 5596.102
$ws.Range("K70").Value = 3652.9998
$ws.Range("L70").Value = 16788.306
$ws.Range("M70").Value = -3382.9998
$ws.Range("N70").Value = -17328.306
$ws.Range("H73").Value = 4916.6895
$ws.Range("I73").Value = 1217.6666
$ws.Range("J73").Value = 5596.102
$ws.Range("K73").Value = 3652.9998
$ws.Range("L73").Value = 16788.306
$ws.Range("M73").Value = -2716.9998
$ws.Range("N73").Value = -18660.306
$ws.Range("H74").Value = 5075.364
$ws.Range("I74").Value = 4730.5
$ws.Range("J74").Value = 5995
$ws.Range("K74").Value = 4730.5
$ws.Range("L74").Value = 5995
$ws.Range("M74").Value = -3794.5
$ws.Range("N74").Value = -7867
$ws.Range("H77").Value = 5075.364
$ws.Range("I77").Value = 4730.5
$ws.Range("J77").Value = 5995
$ws.Range("K77").Value = 23652.5
$ws.Range("L77").Value = 29975
$ws.Range("M77").Value = -18972.5
$ws.Range("N77").Value = -39335
$ws.Range("H131").Value = 3775.5264
$ws.Range("J131").Value = 5296.5
$ws.Range("L131").Value = 15889.5
$ws.Range("N131").Value = -25969.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 59063
$ws.Range("I45").Value = 85185.086
$ws.Range("K45").Value = 85185.086
$ws.Range("M45").Value = -84808.086
$ws.Range("H74").Value = 42479.918
$ws.Range("I74").Value = 5648.8647
$ws.Range("J74").Value = 166366.19
$ws.Range("K74").Value = 5648.8647
$ws.Range("L74").Value = 166366.19
$ws.Range("M74").Value = -4774.8647
$ws.Range("N74").Value = -168114.19
$ws.Range("H77").Value = 42479.918
$ws.Range("I77").Value = 5648.8647
$ws.Range("J77").Value = 166366.19
$ws.Range("K77").Value = 28244.3235
$ws.Range("L77").Value = 831830.95
$ws.Range("M77").Value = -23876.3235
$ws.Range("N77").Value = -840566.95
$ws.Range("H97").Value = 18330.625
$ws.Range("I97").Value = 14636.923
$ws.Range("J97").Value = 34336.668
$ws.Range("K97").Value = 14636.923
$ws.Range("L97").Value = 34336.668
$ws.Range("M97").Value = -14140.923
$ws.Range("N97").Value = -35328.668
$ws.Range("H122").Value = 5700266.5
$ws.Range("I122").Value = 6946347
$ws.Range("J122").Value = 3898.8572
$ws.Range("K122").Value = 20839041
$ws.Range("L122").Value = 11696.5716
$ws.Range("M122").Value = -20836591
$ws.Range("N122").Value = -16596.5716

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1951.909
$ws.Range("I105").Value = 2036
$ws.Range("K105").Value = 2036
$ws.Range("M105").Value = -289
$ws.Range("H134").Value = 2628.4082
$ws.Range("I134").Value = 1196.3414
$ws.Range("K134").Value = 3589.0242
$ws.Range("M134").Value = -1054.0242

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13722.247
$ws.Range("I31").Value = 1514.7
$ws.Range("K31").Value = 1514.7
$ws.Range("M31").Value = -1219.7
$ws.Range("H34").Value = 13722.247
$ws.Range("I34").Value = 1514.7
$ws.Range("K34").Value = 1514.7
$ws.Range("M34").Value = -1312.7
$ws.Range("H99").Value = 4251.75
$ws.Range("J99").Value = 4715.5713
$ws.Range("L99").Value = 4715.5713
$ws.Range("N99").Value = -7711.5713
$ws.Range("H107").Value = 333337300
$ws.Range("I107").Value = 5975
$ws.Range("K107").Value = 5975
$ws.Range("M107").Value = -4055
$ws.Range("H126").Value = 4251.75
$ws.Range("J126").Value = 4715.5713
$ws.Range("L126").Value = 14146.7139
$ws.Range("N126").Value = -19086.7139

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4221.4
$ws.Range("I11").Value = 10012.5
$ws.Range("J11").Value = 360.66666
$ws.Range("K11").Value = 30037.5
$ws.Range("L11").Value = 1081.99998
$ws.Range("M11").Value = -29897.5
$ws.Range("N11").Value = -1361.99998
$ws.Range("H33").Value = 7253
$ws.Range("I33").Value = 87.72727
$ws.Range("K33").Value = 526.3636200000001
$ws.Range("M33").Value = -243.3636200000001
$ws.Range("H41").Value = 1833.3334
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 2250
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 6750
$ws.Range("M41").Value = -2662
$ws.Range("N41").Value = -7426
$ws.Range("H80").Value = 1875.875
$ws.Range("I80").Value = 2227
$ws.Range("J80").Value = 1524.75
$ws.Range("K80").Value = 6681
$ws.Range("L80").Value = 4574.25
$ws.Range("M80").Value = -5745
$ws.Range("N80").Value = -6446.25
$ws.Range("H83").Value = 1875.875
$ws.Range("I83").Value = 2227
$ws.Range("J83").Value = 1524.75
$ws.Range("K83").Value = 20043
$ws.Range("L83").Value = 13722.75
$ws.Range("M83").Value = -15363
$ws.Range("N83").Value = -23082.75
$ws.Range("H107").Value = 1349.25
$ws.Range("J107").Value = 2498.75
$ws.Range("L107").Value = 7496.25
$ws.Range("N107").Value = -11336.25
$ws.Range("H116").Value = 5794.231
$ws.Range("I116").Value = 1573.5
$ws.Range("J116").Value = 7670.1113
$ws.Range("K116").Value = 4720.5
$ws.Range("L116").Value = 23010.3339
$ws.Range("M116").Value = -1278.5
$ws.Range("N116").Value = -29894.3339
$ws.Range("H132").Value = 1690.0416
$ws.Range("I132").Value = 985.0909
$ws.Range("J132").Value = 2286.5386
$ws.Range("K132").Value = 8865.8181
$ws.Range("L132").Value = 20578.8474
$ws.Range("M132").Value = -6335.8181
$ws.Range("N132").Value = -25638.8474
$ws.Range("H134").Value = 2630
$ws.Range("I134").Value = 2630
$ws.Range("K134").Value = 7890
$ws.Range("M134").Value = -2820

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 39996
$ws.Range("J63").Value = 39996
$ws.Range("L63").Value = 39996
$ws.Range("N63").Value = -41368
$ws.Range("H66").Value = 39996
$ws.Range("J66").Value = 39996
$ws.Range("L66").Value = 119988
$ws.Range("N66").Value = -126852
$ws.Range("H97").Value = 3252.5
$ws.Range("I97").Value = 3252.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3252.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2756.5
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 23230.334
$ws.Range("I102").Value = 1893.6945
$ws.Range("K102").Value = 1893.6945
$ws.Range("M102").Value = -271.6945000000001
$ws.Range("H113").Value = 1706.6
$ws.Range("I113").Value = 1674
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1674
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 496
$ws.Range("N113").Value = -6340
$ws.Range("H126").Value = 3285.5386
$ws.Range("I126").Value = 3170.8
$ws.Range("J126").Value = 3668
$ws.Range("K126").Value = 9512.400000000001
$ws.Range("L126").Value = 11004
$ws.Range("M126").Value = -7042.400000000001
$ws.Range("N126").Value = -15944

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2899.8276
$ws.Range("I126").Value = 2766.7083
$ws.Range("K126").Value = 8300.124899999999
$ws.Range("M126").Value = -5830.124899999999
$ws.Range("H136").Value = 2187.3447
$ws.Range("I136").Value = 2022.4584
$ws.Range("J136").Value = 2978.8
$ws.Range("K136").Value = 6067.3752
$ws.Range("L136").Value = 8936.400000000001
$ws.Range("M136").Value = -3517.3752
$ws.Range("N136").Value = -14036.4
